{"js": "const body = context.document.body;\nconst searchResults = body.search(\"Group2\", { matchCase: true, matchWholeWord: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"Group0\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$targetText = \"Group2\"\n$newText = \"Group0\"\n\nforeach ($tbl in $d.Tables) {\n  foreach ($row in $tbl.Rows) {\n    foreach ($cell in $row.Cells) {\n      $cellText = $cell.Range.Text.TrimEnd([char]7, [char]13)\n      if ($cellText -eq $targetText) {\n        $cell.Range.Text = $newText\n      }\n    }\n  }\n}\n"}
